$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.063.19"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.885.51"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.42"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3917"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08417"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.26"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.208"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.27"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.50"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.285"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.04"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001100"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06713"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.953"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.027.15"
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.086.87"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.76"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.36"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.77"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.804"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.657"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02457"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06563"
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.034"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2167"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.083"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.240"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6405"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.19"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6002"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.12"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.673"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.014"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.42"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.165"
$ws.Range("E51").Value = "  -3.29%  "
